# Updated cryptos list on Sat Jun  8 04:08:01 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row with
# the latest scraped values, and re-applies the current ranking order for
# the two pairs of rows whose relative rank flipped since the last run
# (EthereumClassic/Kaspa at rows 34-35, and InjectiveProtocol/dogwifhat at
# rows 47-48), carrying their Coin name/Link along with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells that receive new values to be stored as text,
# matching the original inlineStr type (prevents Excel from auto-converting
# numeric-looking price strings like "1.00" or "681.30" into numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.414.12'
$ws.Range("E2").Value = '  -2.33%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.688.14'
$ws.Range("E3").Value = '  -3.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '681.30'
$ws.Range("E5").Value = '  -3.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.58'
$ws.Range("E6").Value = '  -4.67%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.687.30'
$ws.Range("E7").Value = '  -3.13%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  -4.36%  '

$ws.Range("E10").Value = '  -8.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.34'
$ws.Range("E11").Value = '  -4.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.447'
$ws.Range("E12").Value = '  -2.53%  '

$ws.Range("E13").Value = '  -5.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.60'
$ws.Range("E14").Value = '  -5.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.312.17'
$ws.Range("E15").Value = '  -3.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.686.04'
$ws.Range("E16").Value = '  -3.34%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.451.17'
$ws.Range("E17").Value = '  -2.13%  '

$ws.Range("E18").Value = '  -1.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.32'
$ws.Range("E19").Value = '  -6.32%  '

$ws.Range("E20").Value = '  -6.69%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '482.58'
$ws.Range("E21").Value = '  -3.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.90'
$ws.Range("E22").Value = '  -7.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.667'
$ws.Range("E23").Value = '  -7.71%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.30'
$ws.Range("E24").Value = '  -4.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.834.67'
$ws.Range("E25").Value = '  -3.11%  '

$ws.Range("E26").Value = '  -8.79%  '

$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.49'
$ws.Range("E28").Value = '  -4.60%  '

$ws.Range("E29").Value = '  -7.32%  '

$ws.Range("E30").Value = '  -8.86%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.71'
$ws.Range("E31").Value = '  -10.58%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.88'
$ws.Range("E32").Value = '  -6.06%  '

$ws.Range("E33").Value = '  -7.46%  '

$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.169'
$ws.Range("E34").Value = '  -1.66%  '

$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.13'
$ws.Range("E35").Value = '  -6.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.659.38'
$ws.Range("E37").Value = '  -3.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.50'
$ws.Range("E38").Value = '  -5.98%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.31'
$ws.Range("E39").Value = '  +6.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0939'
$ws.Range("E40").Value = '  -7.01%  '

$ws.Range("E41").Value = '  -4.49%  '

$ws.Range("E43").Value = '  +0.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.955'
$ws.Range("E44").Value = '  -6.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '162.09'
$ws.Range("E45").Value = '  -2.92%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '48.44'
$ws.Range("E46").Value = '  -1.21%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '30.31'
$ws.Range("E47").Value = '  +7.96%  '

$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.84'
$ws.Range("E48").Value = '  -12.73%  '

$ws.Range("E49").Value = '  -7.87%  '

$ws.Range("E50").Value = '  +0.19%  '

$ws.Range("E51").Value = '  -0.92%  '

Write-Host "cryptos list updated"
